$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.817.44"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "2.907.38"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'527.27"
$ws.Range("E5").Value = "  -2.76%  "
$ws.Range("D6").Value = "'144.58"
$ws.Range("E6").Value = "  -4.72%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -3.92%  "
$ws.Range("D9").Value = "2.912.71"
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("E10").Value = "  -5.15%  "
$ws.Range("D11").Value = "'6.13"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "'0.359"
$ws.Range("E12").Value = "  -2.72%  "
$ws.Range("D13").Value = "3.409.21"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("E14").Value = "  +2.69%  "
$ws.Range("D15").Value = "60.843.28"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("D16").Value = "'22.56"
$ws.Range("E16").Value = "  -5.96%  "
$ws.Range("D17").Value = "2.912.48"
$ws.Range("E17").Value = "  -2.85%  "
$ws.Range("E18").Value = "  -3.75%  "
$ws.Range("D19").Value = "'4.90"
$ws.Range("E19").Value = "  -5.01%  "
$ws.Range("D20").Value = "'11.56"
$ws.Range("E20").Value = "  -4.36%  "
$ws.Range("D21").Value = "'354.18"
$ws.Range("E21").Value = "  -6.53%  "
$ws.Range("D22").Value = "'6.54"
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'5.72"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "'64.87"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("E26").Value = "  -3.98%  "
$ws.Range("E27").Value = "  -4.12%  "
$ws.Range("D28").Value = "'1.01"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("E29").Value = "  -4.97%  "
$ws.Range("E30").Value = "  -7.21%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").Value = "'19.66"
$ws.Range("D34").Value = "'153.10"
$ws.Range("E34").Value = "  -4.97%  "
$ws.Range("D35").Value = "'4.41"
$ws.Range("E35").Value = "  -4.15%  "
$ws.Range("D36").Value = "'5.58"
$ws.Range("E36").Value = "  -6.65%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -6.45%  "
$ws.Range("E38").Value = "  -6.18%  "
$ws.Range("D39").Value = "'37.58"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("E40").Value = "  -5.26%  "
$ws.Range("E41").Value = "  -2.90%  "
$ws.Range("E42").Value = "  -5.21%  "
$ws.Range("D43").Value = "2.287.46"
$ws.Range("E43").Value = "  -5.42%  "
$ws.Range("D44").Value = "'0.0582"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("D45").Value = "'20.34"
$ws.Range("E45").Value = "  -7.87%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'4.93"
$ws.Range("E47").Value = "  -4.83%  "
$ws.Range("E48").Value = "  -3.23%  "
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("E50").Value = "  -3.97%  "
$ws.Range("D51").Value = "'18.57"
$ws.Range("E51").Value = "  -6.03%  "
